$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '43.893.26'
Set-TextValue $ws.Range('E2') '  +5.23%  '
Set-TextValue $ws.Range('D3') '2.296.63'
Set-TextValue $ws.Range('E3') '  +3.19%  '
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '232.20'
Set-TextValue $ws.Range('E5') '  +0.42%  '
Set-TextValue $ws.Range('E6') '  +0.48%  '
Set-TextValue $ws.Range('D7') '62.19'
Set-TextValue $ws.Range('E7') '  +2.35%  '
Set-TextValue $ws.Range('E8') '  -0.02%  '
Set-TextValue $ws.Range('D9') '0.418'
Set-TextValue $ws.Range('E9') '  +4.10%  '
Set-TextValue $ws.Range('D10') '0.0920'
Set-TextValue $ws.Range('E10') '  +3.49%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D11') '0.104'
Set-TextValue $ws.Range('E11') '  +0.94%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D12') '2.632.90'
Set-TextValue $ws.Range('E12') '  +3.01%  '
Set-TextValue $ws.Range('D13') '15.88'
Set-TextValue $ws.Range('E13') '  +1.26%  '
Set-TextValue $ws.Range('D14') '24.19'
Set-TextValue $ws.Range('E14') '  +11.31%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D15') '5.75'
Set-TextValue $ws.Range('E15') '  +3.24%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D16') '0.815'
Set-TextValue $ws.Range('E16') '  +1.90%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D17') '2.301.36'
Set-TextValue $ws.Range('E17') '  +3.40%  '
Set-TextValue $ws.Range('D18') '43.662.80'
Set-TextValue $ws.Range('E18') '  +4.84%  '
Set-TextValue $ws.Range('D19') '0.0₃0932'
Set-TextValue $ws.Range('E19') '  +4.53%  '
Set-TextValue $ws.Range('D20') '73.44'
Set-TextValue $ws.Range('E20') '  +0.91%  '
Set-TextValue $ws.Range('E21') '  +3.87%  '
Set-TextValue $ws.Range('D22') '250.62'
Set-TextValue $ws.Range('E22') '  +0.39%  '
Set-TextValue $ws.Range('E24') '  +7.33%  '
Set-TextValue $ws.Range('D25') '2.37'
Set-TextValue $ws.Range('E25') '  +2.59%  '
Set-TextValue $ws.Range('D26') '9.90'
Set-TextValue $ws.Range('E26') '  +3.82%  '
Set-TextValue $ws.Range('D27') '169.82'
Set-TextValue $ws.Range('E27') '  +1.27%  '
Set-TextValue $ws.Range('E28') '  +0.06%  '
Set-TextValue $ws.Range('D29') '20.65'
Set-TextValue $ws.Range('E29') '  +3.59%  '
Set-TextValue $ws.Range('D30') '1.49'
Set-TextValue $ws.Range('E30') '  +5.61%  '
Set-TextValue $ws.Range('E31') '  +0.85%  '
Set-TextValue $ws.Range('E32') '  +0.01%  '
Set-TextValue $ws.Range('D33') '5.02'
Set-TextValue $ws.Range('E33') '  +1.27%  '
Set-TextValue $ws.Range('E34') '  +2.74%  '
Set-TextValue $ws.Range('D35') '0.0660'
Set-TextValue $ws.Range('E35') '  +5.78%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D36') '2.45'
Set-TextValue $ws.Range('E36') '  +3.92%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range('D37') '6.56'
Set-TextValue $ws.Range('E37') '  -1.16%  '
Set-TextValue $ws.Range('D38') '3.67'
Set-TextValue $ws.Range('E38') '  -0.21%  '
Set-TextValue $ws.Range('D39') '0.0252'
Set-TextValue $ws.Range('E39') '  +4.70%  '
Set-TextValue $ws.Range('E40') '  -0.01%  '
Set-TextValue $ws.Range('D41') '8.78'
Set-TextValue $ws.Range('E41') '  +1.53%  '
Set-TextValue $ws.Range('D42') '4.60'
Set-TextValue $ws.Range('E42') '  -4.19%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D43') '0.0975'
Set-TextValue $ws.Range('E43') '  -0.46%  '
$ws.Range('B44').Value = 'TerraClassic'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue $ws.Range('D44') '0.000217'
Set-TextValue $ws.Range('E44') '  -15.18%  '
Set-TextValue $ws.Range('D45') '1.22'
Set-TextValue $ws.Range('E45') '  +0.01%  '
Set-TextValue $ws.Range('D46') '98.51'
Set-TextValue $ws.Range('E46') '  -0.13%  '
Set-TextValue $ws.Range('D47') '1.474.56'
Set-TextValue $ws.Range('E47') '  +0.43%  '
Set-TextValue $ws.Range('D48') '16.75'
Set-TextValue $ws.Range('E48') '  +1.26%  '
Set-TextValue $ws.Range('D49') '2.30'
Set-TextValue $ws.Range('E49') '  +10.30%  '
Set-TextValue $ws.Range('E50') '  +1.75%  '
